# Final Patch Update before release
# Mirror the text from column A into column D for a specific set of rows.
# This matches the source diff, which adds a <c r="D{row}" t="s"><v>...</v></c>
# cell (same shared-string value as the corresponding A{row} cell) for rows
# 39, 42, 44, 51, 56 and 58, and widens the sheet dimension/row spans from
# "A1:C136" / "1:3" to "A1:D136" / "1:4".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(39, 42, 44, 51, 56, 58)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = $ws.Range("A$r").Value2
}
